$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the "P" (D) and "P_variance" (E) values for rows 2-3,
# and the "v" (Q) values for rows 2-3 per the day-ahead schedule tweak.
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 0.3
$ws.Range("Q2").Value = 4

$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0.3
$ws.Range("Q3").Value = 4

# Move the active selection to G9, matching the saved cursor position.
$ws.Range("G9").Select()
